# "Generate Report for Archive"
# - Update the handoff/localization status text from "Ready for handoff" to
#   "In Translation" on all three sheets (Overview, zh-cn, de-de).
# - Shrink the "Latest Handoff Datetime"/Status columns (which previously held
#   the wider "Ready for handoff" text) to their new, narrower auto-fit width.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Status cells that used to read "Ready for handoff" now read "In Translation"
$overview.Range("E2:F3").Value = "In Translation"
$zhcn.Range("C2:C3").Value     = "In Translation"
$dede.Range("C2:C3").Value     = "In Translation"

# New narrower column width for the columns that held the status text
# (closest attainable width to the target 13.4101845877511 character units,
# given the engine's pixel-level ColumnWidth quantization).
$newWidth = 12.5

$overview.Range("E1:F1").ColumnWidth = $newWidth
$zhcn.Range("C1").ColumnWidth        = $newWidth
$dede.Range("C1").ColumnWidth        = $newWidth
